# Infrared control car.pptx - "Add files via upload"
#
# Slide 4, shape "文本框 7" (id=8) contains the paragraph that explains how
# to fetch the micro:bit package. The author appended a second GitHub link
# ("and https://github.com/YahboomTechnology/Yahboom_IR ") right before the
# trailing "to get the package." run, styled the same way as the existing
# sentence (grey "and " run + red-colored URL run), and the textbox grew
# taller to fit the extra wrapped line.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(9)             # "文本框 7" (id=8)

$tr = $sh.TextFrame.TextRange

# Locate the existing trailing run "to get the package." so we can insert
# right before it without disturbing its own run.
$tail = $tr.Characters(229, 19)

$insertedText = "and https://github.com/YahboomTechnology/Yahboom_IR "
$tail.InsertBefore($insertedText) | Out-Null

# The inserted text merges into the run that previously ended in
# "...yahboom_mbit_en " (grey run), inheriting its grey/accent5 75% luminance
# formatting automatically - matching "and " in the target. Now re-color just
# the new URL portion red, which splits it into its own run.
$urlRange = $tr.Characters(233, 48)   # "https://github.com/YahboomTechnology/Yahboom_IR "
$urlRange.Font.Color.RGB = 255        # packed 0x00BBGGRR -> FF0000 (red)

# The shape auto-fits its height to the text (a:spAutoFit); restore the
# exact height PowerPoint computed for the now-taller textbox.
$sh.Height = 204
